# The "#export" sheet's row 14 held a "#sample%untrack=project.id" header
# (shared-string index 16) describing an untrack directive with only the
# project.id field. The edit extends it to also untrack the generated
# project.id%number field, matching the sibling "#sample%track=..." header
# directly above it (row 1), which already lists both fields.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = "#sample%untrack=project.id,project.id%number"

# Move/restore the worksheet's active selection to the cell that was edited.
$ws.Range("C14").Select()
